# Applies the inventory update + window/layout cleanup described in the commit message:
# "Added updating quantities, and organized inventory update window"

$wb = $excel.ActiveWorkbook

$wsInventory = $wb.Worksheets.Item("Sheet1")
$wsAmazon = $wb.Worksheets.Item("Daily Amazon")

# --- Sheet1 ("Inventory") updates ---

# Update quantity for Juzo (row 8)
$wsInventory.Range("C8").Value = 30

# Delete the duplicate "Isabela" / "Naruto with Rasengan" rows (16 and 17)
$wsInventory.Rows.Item(16).Delete()
$wsInventory.Rows.Item(16).Delete()

# Delete the stray "goku" row which is now row 23 after the above deletions
$wsInventory.Rows.Item(23).Delete()

# Select / activate the now-empty row under the data as the new entry point
$wsInventory.Activate()
$wsInventory.Range("A23:C23").Select()

# --- Daily Amazon sheet: scroll position / active cell cleanup ---

$wsAmazon.Activate()
$wsAmazon.Application.ActiveWindow.ScrollRow = 34
$wsAmazon.Range("A66").Select()

# Fix the style of the last few date cells (rows 65-67, column A) so they use the
# standard date-format style (numFmt 164, same as the rest of the Date column)
# instead of the now-removed duplicate date format (numFmt 165)
$wsAmazon.Range("A65:A67").NumberFormat = "yyyy\-mm\-dd"

# --- Workbook window arrangement ---
$wsInventory.Activate()
$excel.ActiveWindow.WindowState = -4137
$excel.Windows.Item(1).WindowState = -4143
